$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header row, bold style matching W1/X1) ---
$ws.Range("Y1").Value = ":8-2015"
$ws.Range("Y1").Font.Bold = $true

# --- Row 2 (header row, bold style matching W2/X2) ---
$ws.Range("Y2").Value = ":8-2015"
$ws.Range("Y2").Font.Bold = $true
$ws.Range("Z2").Value = ":9-2015"
$ws.Range("Z2").Font.Bold = $true

# --- Row 3 (data row) ---
$ws.Range("Y3").Value = 0.442
$ws.Range("Z3").Value = 0.558

# --- Row 4 (data row) ---
$ws.Range("Y4").Value = 0.442
$ws.Range("Z4").Value = 0.558

# --- Row 5 (data row) ---
$ws.Range("Y5").Value = 0
$ws.Range("Z5").Value = 0

# --- Row 6 (data row) ---
$ws.Range("Y6").Value = "DI"
$ws.Range("Z6").Value = "DI"

# Match the final selection recorded in the saved workbook
$ws.Range("Z2").Select()
